$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "27.061.05"
Set-TextValue "E2" "  -3.10%  "
Set-TextValue "D3" "1.715.46"
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "310.26"
Set-TextValue "E5" "  -5.62%  "
Set-TextValue "E6" "  -0.01%  "
Set-TextValue "D7" "0.4797"
Set-TextValue "E7" "  +6.09%  "
Set-TextValue "D8" "0.3451"
Set-TextValue "E8" "  -2.21%  "
Set-TextValue "D9" "42.17"
Set-TextValue "E9" "  +0.61%  "
Set-TextValue "D10" "0.07252"
Set-TextValue "E10" "  -1.76%  "
Set-TextValue "D11" "1.038"
Set-TextValue "E11" "  -5.18%  "
Set-TextValue "E12" "  +0.01%  "
Set-TextValue "D13" "19.73"
Set-TextValue "E13" "  -4.81%  "
Set-TextValue "D14" "5.836"
Set-TextValue "E14" "  -2.93%  "
Set-TextValue "D15" "1.713.96"
Set-TextValue "E15" "  -3.19%  "
Set-TextValue "D16" "6.824"
Set-TextValue "E16" "  -5.13%  "
Set-TextValue "D17" "87.17"
Set-TextValue "E17" "  -5.95%  "
Set-TextValue "E18" "  -2.48%  "
Set-TextValue "D19" "0.06374"
Set-TextValue "E19" "  -0.83%  "
Set-TextValue "E20" "  -0.03%  "
Set-TextValue "D22" "5.621"
Set-TextValue "E22" "  -2.50%  "
Set-TextValue "D23" "27.117.29"
Set-TextValue "E23" "  -2.97%  "
Set-TextValue "E24" "  -4.17%  "
Set-TextValue "D25" "2.099"
Set-TextValue "E25" "  +0.06%  "
Set-TextValue "D26" "19.99"
Set-TextValue "E26" "  -0.95%  "
Set-TextValue "D27" "150.73"
Set-TextValue "E27" "  -5.68%  "
Set-TextValue "D28" "1.910.30"
Set-TextValue "D29" "2.055"
Set-TextValue "E29" "  -3.49%  "
Set-TextValue "E30" "  -2.86%  "
Set-TextValue "D31" "1.035"
Set-TextValue "E31" "  -4.38%  "
Set-TextValue "D32" "0.09251"
Set-TextValue "E32" "  +0.48%  "
Set-TextValue "E33" "  -1.93%  "
Set-TextValue "D34" "5.306"
Set-TextValue "E34" "  -5.36%  "
Set-TextValue "D35" "1.477"
Set-TextValue "E35" "  +6.82%  "
Set-TextValue "D36" "0.02177"
Set-TextValue "E36" "  -4.31%  "
Set-TextValue "D37" "0.05833"
Set-TextValue "E37" "  -4.28%  "
Set-TextValue "D38" "10.91"
Set-TextValue "E38" "  -7.70%  "
Set-TextValue "E39" "  -4.96%  "
Set-TextValue "E40" "  -0.04%  "
Set-TextValue "D41" "4.702"
Set-TextValue "E41" "  -4.98%  "
Set-TextValue "D42" "0.5939"
Set-TextValue "E42" "  -4.92%  "
Set-TextValue "D43" "1.083"
Set-TextValue "E43" "  -8.05%  "
Set-TextValue "D44" "7.486"
Set-TextValue "E44" "  -4.03%  "
Set-TextValue "D45" "12.75"
Set-TextValue "E45" "  -3.59%  "
Set-TextValue "D46" "3.584"
Set-TextValue "E46" "  -4.02%  "
Set-TextValue "D47" "0.5557"
Set-TextValue "E47" "  -4.75%  "
Set-TextValue "D48" "118.49"
Set-TextValue "E48" "  -3.31%  "
Set-TextValue "D49" "1.822"
Set-TextValue "E49" "  -5.46%  "
Set-TextValue "D50" "0.06633"
Set-TextValue "E50" "  -2.97%  "
Set-TextValue "D51" "1.088"
Set-TextValue "E51" "  -4.17%  "
